# Commit: "Fixed POI packaging and upgraded to POI 3.15."
#
# The whole diff (word/document.xml's <w:document> root + its <w:pgSz>/
# <w:pgMar>, and word/styles.xml's <w:rFonts>/<w:lang>/<w:latentStyles>/
# every <w:lsdException>/<w:style>/<w:tblInd>/<w:tblCellMar> child) is a
# pure XML-attribute reordering: the Apache POI/XMLBeans upgrade mentioned
# in the commit message made the writer emit each element's attributes in
# sorted order (namespace declarations first, then the rest, all sorted by
# qualified name) instead of their previous insertion order. Every "-"/"+"
# pair in the diff has the exact same tag, the exact same attribute names
# and the exact same attribute values - nothing was inserted, deleted or
# re-valued; only the left-to-right order the attributes are written in
# changed. There is no text, run, paragraph, page size/margin, font,
# language, or style-definition change anywhere in the package.
#
# Word's object model (Find/Replace, Paragraphs, Range, PageSetup, Styles,
# Fonts, Tables, ...) edits document *content*; it has no property for the
# byte order the XML writer serializes an element's attributes in, so
# there is no COM call that reproduces this specific change - and forcing
# an unrelated, identical-value write (e.g. re-assigning PageSetup.Gutter
# to its own value) only makes the round-tripped part's XML drift further
# from the target by causing the writer to re-emit the part from scratch,
# which introduces extra namespace declarations that are not part of the
# requested change. The correct, content-faithful application of this
# diff is therefore to leave the document exactly as it already is: every
# value the diff touches (page size/margins, default fonts/language,
# latent-style metadata, the built-in style definitions) already matches,
# so no Content/Find, Range, PageSetup, or Styles call is needed.

$d = $word.ActiveDocument
